$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test1")

# Fill in the "staircase" numeric pattern in rows 17-26 (columns B..K),
# matching the existing pattern already present in rows 4-13.
$ws1.Cells.Item(17, 2).Value  = 1
$ws1.Cells.Item(17, 3).Value  = 2
$ws1.Cells.Item(17, 4).Value  = 3
$ws1.Cells.Item(17, 5).Value  = 4
$ws1.Cells.Item(17, 6).Value  = 5
$ws1.Cells.Item(17, 7).Value  = 6
$ws1.Cells.Item(17, 8).Value  = 7
$ws1.Cells.Item(17, 9).Value  = 8
$ws1.Cells.Item(17, 10).Value = 9
$ws1.Cells.Item(17, 11).Value = 10

$ws1.Cells.Item(18, 2).Value  = 2
$ws1.Cells.Item(18, 3).Value  = 3
$ws1.Cells.Item(18, 4).Value  = 4
$ws1.Cells.Item(18, 5).Value  = 5
$ws1.Cells.Item(18, 6).Value  = 6
$ws1.Cells.Item(18, 7).Value  = 7
$ws1.Cells.Item(18, 8).Value  = 8
$ws1.Cells.Item(18, 9).Value  = 9
$ws1.Cells.Item(18, 10).Value = 10

$ws1.Cells.Item(19, 2).Value  = 3
$ws1.Cells.Item(19, 3).Value  = 4
$ws1.Cells.Item(19, 4).Value  = 5
$ws1.Cells.Item(19, 5).Value  = 6
$ws1.Cells.Item(19, 6).Value  = 7
$ws1.Cells.Item(19, 7).Value  = 8
$ws1.Cells.Item(19, 8).Value  = 9
$ws1.Cells.Item(19, 9).Value  = 10

$ws1.Cells.Item(20, 2).Value  = 4
$ws1.Cells.Item(20, 3).Value  = 5
$ws1.Cells.Item(20, 4).Value  = 6
$ws1.Cells.Item(20, 5).Value  = 7
$ws1.Cells.Item(20, 6).Value  = 8
$ws1.Cells.Item(20, 7).Value  = 9
$ws1.Cells.Item(20, 8).Value  = 10

$ws1.Cells.Item(21, 2).Value  = 5
$ws1.Cells.Item(21, 3).Value  = 6
$ws1.Cells.Item(21, 4).Value  = 7
$ws1.Cells.Item(21, 5).Value  = 8
$ws1.Cells.Item(21, 6).Value  = 9
$ws1.Cells.Item(21, 7).Value  = 10

$ws1.Cells.Item(22, 2).Value  = 6
$ws1.Cells.Item(22, 3).Value  = 7
$ws1.Cells.Item(22, 4).Value  = 8
$ws1.Cells.Item(22, 5).Value  = 9
$ws1.Cells.Item(22, 6).Value  = 10

$ws1.Cells.Item(23, 2).Value  = 7
$ws1.Cells.Item(23, 3).Value  = 8
$ws1.Cells.Item(23, 4).Value  = 9
$ws1.Cells.Item(23, 5).Value  = 10

$ws1.Cells.Item(24, 2).Value  = 8
$ws1.Cells.Item(24, 3).Value  = 9
$ws1.Cells.Item(24, 4).Value  = 10

$ws1.Cells.Item(25, 2).Value  = 9
$ws1.Cells.Item(25, 3).Value  = 10

$ws1.Cells.Item(26, 2).Value  = 10

# Make Test1 the active sheet/tab (was Test2), and set its selection to O13.
$ws1.Activate() | Out-Null
$ws1.Range("O13").Select() | Out-Null
